$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": update the two timeslice-group strings (C13/C14).
# G7 (=C14) and G8 (=C13) recalc automatically from these.
$wsUC = $wb.Worksheets.Item("ev_charging_uc")
$wsUC.Range("C13").Value = "FaP,SaP,RaD,RaP,WaD,FaD,SaD,WaP"
$wsUC.Range("C14").Value = "RaP,WaP,FaP,SaP,FaN,RaN,SaN,WaN"

# --- Sheet "re_profiles": re-shuffle the season rows in M4:N7 (hydro block).
$wsRE = $wb.Worksheets.Item("re_profiles")
$wsRE.Range("M4").Value = "W"
$wsRE.Range("N4").Value = 0.22555529847292924
$wsRE.Range("M5").Value = "S"
$wsRE.Range("N5").Value = 0.40439611291068944
$wsRE.Range("M6").Value = "F"
$wsRE.Range("N6").Value = 0.26702915316982878
$wsRE.Range("M7").Value = "R"
$wsRE.Range("N7").Value = 0.30301943544655252
